$wb = $excel.ActiveWorkbook

# --- Sheet 1: status --- add row 3 (duplicate of row 2)
$ws = $wb.Worksheets.Item(1)
$ws.Range("A3").Value = 1465993370
$ws.Range("B3").Value = 1041368833

# --- Sheet 2: neighbors --- add rows 4 and 5 (duplicates of rows 2 and 3)
$ws = $wb.Worksheets.Item(2)
$ws.Range("A4").Value = 1465993370
$ws.Range("B4").Value = "10.0.0.5"
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 1

$ws.Range("A5").Value = 1465993370
$ws.Range("B5").Value = "10.0.0.3"
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 0

# --- Sheet 3: links --- add rows 4 and 5 (duplicates of rows 2 and 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("A4").Value = 1465993370
$ws.Range("B4").Value = "10.0.0.4"
$ws.Range("C4").Value = "10.0.0.5"
$ws.Range("D4").Value = 39683
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1024

$ws.Range("A5").Value = 1465993370
$ws.Range("B5").Value = "10.0.0.4"
$ws.Range("C5").Value = "10.0.0.3"
$ws.Range("D5").Value = 38439
$ws.Range("E5").Value = 0.419
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2440

# --- Sheet 4: routes --- add rows 5, 6, 7 (duplicates of rows 2, 3, 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("A5").Value = 1465993370
$ws.Range("B5").Value = "10.0.0.3"
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = "10.0.0.3"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 2440
$ws.Range("G5").Value = "mesh0"

$ws.Range("A6").Value = 1465993370
$ws.Range("B6").Value = "10.0.0.5"
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = "10.0.0.5"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1024
$ws.Range("G6").Value = "mesh0"

$ws.Range("A7").Value = 1465993370
$ws.Range("B7").Value = "10.0.0.6"
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = "10.0.0.5"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2327
$ws.Range("G7").Value = "mesh0"

# --- Sheet 5: topology --- add rows 8-13 (duplicates of rows 2-7)
$ws = $wb.Worksheets.Item(5)
$ws.Range("A8").Value = 1465993370
$ws.Range("B8").Value = "10.0.0.4"
$ws.Range("C8").Value = "10.0.0.3"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.259
$ws.Range("F8").Value = 3956
$ws.Range("G8").Value = 280425

$ws.Range("A9").Value = 1465993370
$ws.Range("B9").Value = "10.0.0.3"
$ws.Range("C9").Value = "10.0.0.4"
$ws.Range("D9").Value = 0.419
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 2440
$ws.Range("G9").Value = 0

$ws.Range("A10").Value = 1465993370
$ws.Range("B10").Value = "10.0.0.5"
$ws.Range("C10").Value = "10.0.0.4"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1024
$ws.Range("G10").Value = 0

$ws.Range("A11").Value = 1465993370
$ws.Range("B11").Value = "10.0.0.4"
$ws.Range("C11").Value = "10.0.0.5"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1024
$ws.Range("G11").Value = 279823

$ws.Range("A12").Value = 1465993370
$ws.Range("B12").Value = "10.0.0.6"
$ws.Range("C12").Value = "10.0.0.5"
$ws.Range("D12").Value = 0.886
$ws.Range("E12").Value = 0.886
$ws.Range("F12").Value = 1303
$ws.Range("G12").Value = 279823

$ws.Range("A13").Value = 1465993370
$ws.Range("B13").Value = "10.0.0.5"
$ws.Range("C13").Value = "10.0.0.6"
$ws.Range("D13").Value = 0.886
$ws.Range("E13").Value = 0.886
$ws.Range("F13").Value = 1303
$ws.Range("G13").Value = 277804
